$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
$ws1.Cells.Item(17, 1).NumberFormat = $ws1.Cells.Item(16, 1).NumberFormat
$ws1.Cells.Item(17, 1).Value = 45732.76288380787
$ws1.Cells.Item(17, 2).Value = "0x01,0x90"
$ws1.Cells.Item(17, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws1.Cells.Item(17, 4).Value = "0x01,0x86"
$ws1.Cells.Item(17, 5).Value = "0xd"
$ws1.Cells.Item(17, 6).Value = 400
$ws1.Cells.Item(17, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item(17, 8).Value = 390
$ws1.Cells.Item(17, 9).Value = 13

$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
$ws2.Cells.Item(17, 1).NumberFormat = $ws2.Cells.Item(16, 1).NumberFormat
$ws2.Cells.Item(17, 1).Value = 45732.61599116898
$ws2.Cells.Item(17, 2).Value = "0x01,0x90"
$ws2.Cells.Item(17, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(17, 4).Value = "0x01,0x86"
$ws2.Cells.Item(17, 5).Value = "0xe"
$ws2.Cells.Item(17, 6).Value = 400
$ws2.Cells.Item(17, 7).Value = [double]"5.68631262647114e+23"
$ws2.Cells.Item(17, 8).Value = 390
$ws2.Cells.Item(17, 9).Value = 14

$ws3 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
$ws3.Cells.Item(17, 1).NumberFormat = $ws3.Cells.Item(16, 1).NumberFormat
$ws3.Cells.Item(17, 1).Value = 45732.76374813657
$ws3.Cells.Item(17, 2).Value = "0x01,0x90"
$ws3.Cells.Item(17, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item(17, 4).Value = "0x01,0x86"
$ws3.Cells.Item(17, 5).Value = "0x3"
$ws3.Cells.Item(17, 6).Value = 400
$ws3.Cells.Item(17, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item(17, 8).Value = 390
$ws3.Cells.Item(17, 9).Value = 3

$ws4 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
$ws4.Cells.Item(17, 1).NumberFormat = $ws4.Cells.Item(16, 1).NumberFormat
$ws4.Cells.Item(17, 1).Value = 45732.81838934027
$ws4.Cells.Item(17, 2).Value = "0x01,0x90"
$ws4.Cells.Item(17, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Cells.Item(17, 4).Value = "0x01,0x86"
$ws4.Cells.Item(17, 5).Value = "0x3"
$ws4.Cells.Item(17, 6).Value = 400
$ws4.Cells.Item(17, 7).Value = [double]"9.85046333984776e+23"
$ws4.Cells.Item(17, 8).Value = 390
$ws4.Cells.Item(17, 9).Value = 3
